$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.148.06"
$ws.Range("E2").Value = "  +2.22%  "
$ws.Range("D3").Value = "3.133.76"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  +0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.89"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +3.10%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "3.127.30"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.533"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  +12.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.76"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.50"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.41%  "
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("E16").Value = "  +0.78%  "
$ws.Range("D17").Value = "63.932.02"
$ws.Range("E17").Value = "  +1.97%  "
$ws.Range("E18").Value = "  -1.50%  "
$ws.Range("D19").Value = "3.131.65"
$ws.Range("E19").Value = "  +0.76%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "468.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.23%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.37"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.17%  "
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.27"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.48"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.01"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +9.42%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.71"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.78%  "
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.88"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("E33").Value = "  -2.47%  "
$ws.Range("D34").Value = "0.0₃0895"
$ws.Range("E34").Value = "  +11.62%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.39"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +7.63%  "
$ws.Range("E36").Value = "  +1.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.41"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +12.07%  "
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.96"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "456.32"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +7.13%  "
$ws.Range("E41").Value = "  -0.97%  "
$ws.Range("E42").Value = "  +0.34%  "
$ws.Range("D43").Value = "2.898.15"
$ws.Range("E43").Value = "  -1.09%  "
$ws.Range("E44").Value = "  -1.00%  "
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "36.03"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.49"
$ws.Range("D48").ClearFormats()
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E51").Value = "  +0.70%  "
